# Update "想去人数" (interest count) values in column F across sheets,
# matching the data refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 206
$ws1.Range("F6").Value  = 24
$ws1.Range("F7").Value  = 1056
$ws1.Range("F8").Value  = 843
$ws1.Range("F9").Value  = 251
$ws1.Range("F12").Value = 843
$ws1.Range("F13").Value = 290
$ws1.Range("F14").Value = 587
$ws1.Range("F16").Value = 1337
$ws1.Range("F18").Value = 1261
$ws1.Range("F19").Value = 1199
$ws1.Range("F20").Value = 2892
$ws1.Range("F21").Value = 1445
$ws1.Range("F22").Value = 709
$ws1.Range("F26").Value = 1022
$ws1.Range("F27").Value = 362
$ws1.Range("F28").Value = 3143

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 15
$ws2.Range("F6").Value = 11

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 748

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 748
$ws4.Range("F7").Value  = 15
$ws4.Range("F10").Value = 206
$ws4.Range("F11").Value = 11
$ws4.Range("F12").Value = 24
$ws4.Range("F13").Value = 1056
$ws4.Range("F14").Value = 843
$ws4.Range("F15").Value = 251
$ws4.Range("F23").Value = 843
$ws4.Range("F24").Value = 290
$ws4.Range("F25").Value = 587
$ws4.Range("F27").Value = 1337
$ws4.Range("F29").Value = 1261
$ws4.Range("F30").Value = 1199
$ws4.Range("F31").Value = 2892
$ws4.Range("F32").Value = 1445
$ws4.Range("F33").Value = 709
$ws4.Range("F39").Value = 1022
$ws4.Range("F40").Value = 362
$ws4.Range("F41").Value = 3143
